# Regression-model test routine: revert the experimental UK regression
# estimates that had been added to the EUROMOD policy schedule, keeping
# only the two original rows (2015, 2019) on the "UK" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")

# Sheet currently holds header (row 1) + 17 data rows (rows 2-18, years
# 2011-2027). Only the rows for 2015 (currently row 6) and 2019
# (currently row 10) should remain, ending up as the new rows 2 and 3.
# Delete the unwanted rows bottom-to-top so row numbers of the
# not-yet-deleted rows stay stable while iterating.
$rowsToDelete = @(18, 17, 16, 15, 14, 13, 12, 11, 9, 8, 7, 5, 4, 3, 2)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

Write-Host "UK sheet now has" $ws.UsedRange.Rows.Count "rows"
